# Patient data template update:
#  - Replace sample row values (row 2) with a new, non-identifying sample
#  - Clear several demographic columns in the sample row
#  - Drop the custom column widths so columns use the default width

$wb = $excel.ActiveWorkbook

# Remember original sheet name/position before adding a helper sheet.
$old = $wb.Worksheets.Item(1)
$oldName = $old.Name

# Add a brand-new worksheet. New sheets never carry over the custom
# <cols> widths that were baked into the original sheet, so building the
# replacement data here and then dropping the old sheet is the simplest
# way to shed those widths.
$new = $wb.Worksheets.Add()

# After the Add(), the brand-new (blank) sheet sits at index 1 and the
# original sheet has been pushed to index 2. Resolve both "live" from
# the collection (by position) instead of relying on stale references.
$src = $wb.Worksheets.Item(2)
$dst = $wb.Worksheets.Item(1)

# Copy the header row (row 1) through unchanged.
for ($col = 1; $col -le 13; $col++) {
    $dst.Cells.Item(1, $col).Value = $src.Cells.Item(1, $col).Value2
}

# Row 2: new sample patient values.
$dst.Cells.Item(2, 1).Value  = "Maria Silva"                     # A2 nome
$dst.Cells.Item(2, 2).Value  = "QR 100 Conjunto 5 Casa 20"        # B2 endereco
$dst.Cells.Item(2, 3).Value  = "72000-000"                        # C2 cep
$dst.Cells.Item(2, 4).Value  = ""                                 # D2 nome_social
$dst.Cells.Item(2, 5).Value  = ""                                 # E2 nome_mae
$dst.Cells.Item(2, 6).Value  = ""                                 # F2 data_nascimento
$dst.Cells.Item(2, 7).Value  = ""                                 # G2 idade
$dst.Cells.Item(2, 8).Value  = ""                                 # H2 cns_ou_cpf
$dst.Cells.Item(2, 9).Value  = "(61) 98888-8888"                  # I2 telefone
$dst.Cells.Item(2, 10).Value = ""                                 # J2 identidade_genero
$dst.Cells.Item(2, 11).Value = ""                                 # K2 cor_raca
$dst.Cells.Item(2, 12).Value = ""                                 # L2 condicoes_saude
$dst.Cells.Item(2, 13).Value = ""                                 # M2 ultimo_atendimento

# Remove the old sheet (now at index 2) and rename the new one back to
# the original sheet name/position.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item(2).Delete() | Out-Null
$wb.Worksheets.Item(1).Name = $oldName
